$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Starting point: a single paragraph reading "Hello friend ." with the
# (hidden) "_GoBack" bookmark Word leaves at the last edit point, right
# after that text.
#
# Target:
#   Para 1: "Hello " | proofErr(gramStart) | "friend ." | proofErr(gramEnd)
#   Para 2: "Now bed time." followed by the "_GoBack" bookmark.
# ------------------------------------------------------------------

$p1 = $d.Paragraphs(1)

# Remember the original paragraph/text boundary *before* touching the
# document: this is where the new paragraph break belongs (right after
# "Hello friend ." and before "Now bed time.").
$splitPos = $p1.Range.End - 1

# Find the point to work from: the "_GoBack" bookmark if Word left one
# (it marks exactly where the last edit happened), otherwise just the
# end of the first paragraph's text (before its paragraph mark).
if ($d.Bookmarks.Exists("_GoBack")) {
    $editPoint = $d.Bookmarks("_GoBack").Range
} else {
    $editPoint = $d.Range($splitPos, $splitPos)
    $editPoint.Collapse(0)
}

# 1) Type "Now bed time." right at that point. Inserting plain text
#    exactly at a collapsed bookmark pushes the bookmark to trail the
#    freshly typed text (the same gravity real Word uses), so
#    "_GoBack" ends up right after "Now bed time.".
$editPoint.InsertAfter("Now bed time.")

# 2) Press Enter at the original boundary to turn that one paragraph
#    into two: "Hello friend ." / "Now bed time." (+bookmark).
$breakPoint = $d.Range($splitPos, $splitPos)
$breakPoint.InsertParagraphAfter()

# 3) Mark the space-before-period in "friend ." as a grammar slip, the
#    way Word's proofer brackets it: split the single run into
#    "Hello " + proofErr(gramStart) + "friend ." + proofErr(gramEnd).
$p1 = $d.Paragraphs(1)
$runRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$proofedXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Hello </w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>friend .</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$runRange.InsertXML($proofedXml)

Write-Output "Paragraph 1: [$($d.Paragraphs(1).Range.Text)]"
Write-Output "Paragraph 2: [$($d.Paragraphs(2).Range.Text)]"
